$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# New measurement rows (low-frequency oscilloscope readings), added after row 29
# Columns: A=f(kHz) B=Vin pk pk(V) C=Vout pk pk(V) D=H dB (formula) E=phase/deg
$ws.Range("A30").Value = 1000
$ws.Range("B30").Value = 20
$ws.Range("C30").Value = 0.52
$ws.Range("D30").Formula = "=20*LOG(C30/B30)"
$ws.Range("E30").Value = -86

$ws.Range("A31").Value = 500
$ws.Range("B31").Value = 20
$ws.Range("C31").Value = 1.03
$ws.Range("D31").Formula = "=20*LOG(C31/B31)"
$ws.Range("E31").Value = -85

$ws.Range("E31").Select()
